$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.614.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.859.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6946"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07704"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3061"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07780"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.152"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.856.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6929"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.572"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.483.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008307"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.101.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.616"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  +2.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.933"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.535"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.256"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.181"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.202"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05132"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("E35").Value = "  +3.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.155"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.687"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.333.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01874"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.727"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9723"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.806"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000127"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.46%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.785"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.53%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.999.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5219"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.782"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.967"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.01%  "
